{"js": "// Replace the date header and each division prompt in the table with the\n// values from the new practice sheet. Each old string occurs exactly once\n// in the document, so a scoped search + replace per pair is safe and keeps\n// the original run formatting (font/size) intact.\nconst replacements = [\n  [\"2025-11-21 Friday\", \"2025-11-22 Saturday\"],\n  [\"335\u00f78=\", \"431\u00f74=\"],\n  [\"307\u00f77=\", \"572\u00f74=\"],\n  [\"166\u00f74=\", \"437\u00f78=\"],\n  [\"678\u00f77=\", \"408\u00f77=\"],\n  [\"801\u00f72=\", \"382\u00f74=\"],\n  [\"612\u00f76=\", \"957\u00f79=\"],\n  [\"304\u00f75=\", \"508\u00f72=\"],\n  [\"725\u00f73=\", \"253\u00f74=\"],\n  [\"404\u00f75=\", \"577\u00f73=\"],\n  [\"261\u00f72=\", \"234\u00f76=\"],\n  [\"479\u00f74=\", \"462\u00f79=\"],\n  [\"814\u00f73=\", \"878\u00f73=\"],\n  [\"788\u00f76=\", \"154\u00f74=\"],\n  [\"437\u00f77=\", \"926\u00f78=\"],\n  [\"962\u00f78=\", \"883\u00f76=\"],\n  [\"629\u00f72=\", \"514\u00f72=\"],\n  [\"595\u00f76=\", \"523\u00f77=\"],\n  [\"461\u00f76=\", \"453\u00f78=\"],\n  [\"746\u00f75=\", \"454\u00f77=\"],\n  [\"972\u00f74=\", \"297\u00f72=\"],\n  [\"995\u00f73=\", \"815\u00f76=\"],\n  [\"140\u00f77=\", \"755\u00f73=\"],\n  [\"776\u00f73=\", \"956\u00f72=\"],\n  [\"361\u00f78=\", \"950\u00f72=\"],\n  [\"878\u00f74=\", \"770\u00f79=\"],\n];\n\nconst body = context.document.body;\nconst searchResults = [];\n\nfor (const [oldText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  searchResults.push(results);\n}\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searchResults[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date header and each division prompt in the table with the\n# values from the new practice sheet. Each old string occurs exactly once\n# in the document, so Find/Replace per pair is safe and preserves the\n# original run formatting (font/size) of the matched text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-11-21 Friday\", \"2025-11-22 Saturday\"),\n  @(\"335\u00f78=\", \"431\u00f74=\"),\n  @(\"307\u00f77=\", \"572\u00f74=\"),\n  @(\"166\u00f74=\", \"437\u00f78=\"),\n  @(\"678\u00f77=\", \"408\u00f77=\"),\n  @(\"801\u00f72=\", \"382\u00f74=\"),\n  @(\"612\u00f76=\", \"957\u00f79=\"),\n  @(\"304\u00f75=\", \"508\u00f72=\"),\n  @(\"725\u00f73=\", \"253\u00f74=\"),\n  @(\"404\u00f75=\", \"577\u00f73=\"),\n  @(\"261\u00f72=\", \"234\u00f76=\"),\n  @(\"479\u00f74=\", \"462\u00f79=\"),\n  @(\"814\u00f73=\", \"878\u00f73=\"),\n  @(\"788\u00f76=\", \"154\u00f74=\"),\n  @(\"437\u00f77=\", \"926\u00f78=\"),\n  @(\"962\u00f78=\", \"883\u00f76=\"),\n  @(\"629\u00f72=\", \"514\u00f72=\"),\n  @(\"595\u00f76=\", \"523\u00f77=\"),\n  @(\"461\u00f76=\", \"453\u00f78=\"),\n  @(\"746\u00f75=\", \"454\u00f77=\"),\n  @(\"972\u00f74=\", \"297\u00f72=\"),\n  @(\"995\u00f73=\", \"815\u00f76=\"),\n  @(\"140\u00f77=\", \"755\u00f73=\"),\n  @(\"776\u00f73=\", \"956\u00f72=\"),\n  @(\"361\u00f78=\", \"950\u00f72=\"),\n  @(\"878\u00f74=\", \"770\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
